# Add 2022-Q4 data
# 1) Insert a brand-new worksheet right after "总计" (position 2), holding the
#    2022-Q4 per-fund detail (mirrors the layout already used by the other
#    quarterly sheets), and name it "2022-Q4".
# 2) Insert a new row at the top of the "总计" summary sheet for the 2022-Q4
#    totals, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# --- 1. Create the "2022-Q4" worksheet, placed right after "总计" ---------
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item(1))
$newSheet.Name = "2022-Q4"

# Worksheets.Item(...) handles returned above are bound to a *position*, not
# a specific sheet, so re-fetch everything we still need now that the sheet
# count/order has settled.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item("2022-Q4")
$q4_2021Sheet = $wb.Worksheets.Item("2021-Q4")

# Header row (same headers every quarterly sheet uses)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Single fund row for the quarter. The fund code and the percentage-ish
# figures are stored as plain TEXT in the source data (not numbers), so
# force a text number-format before assigning them — otherwise Excel helpfully
# "upgrades" a numeric-looking string to a real number. A blank, never-touched
# cell (Z99) is used afterwards to paste the plain default format back on top,
# so the text stays text without leaving a stray "@" number-format behind.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "003956"
$newSheet.Range("C2").Value = "南方产业智选股票"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "3.60"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "85.80"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "4.89"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.1760"
$newSheet.Range("H2").Value = 4

$newSheet.Range("Z99").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)
$newSheet.Range("Z99").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)
$newSheet.Range("Z99").Copy()
$newSheet.Range("E2").PasteSpecial(-4122)
$newSheet.Range("Z99").Copy()
$newSheet.Range("F2").PasteSpecial(-4122)
$newSheet.Range("Z99").Copy()
$newSheet.Range("G2").PasteSpecial(-4122)
$newSheet.Range("Z99").Clear()

# Match the bold/centered/bordered header style used elsewhere (copy format
# from the equivalent cells on the existing "2021-Q4" sheet).
$q4_2021Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4_2021Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# --- 2. Insert the 2022-Q4 summary row into "总计", pushing others down ---
$b4 = $totalSheet.Range("B4").Value()
$c4 = $totalSheet.Range("C4").Value()
$d4 = $totalSheet.Range("D4").Value()

$b3 = $totalSheet.Range("B3").Value()
$c3 = $totalSheet.Range("C3").Value()
$d3 = $totalSheet.Range("D3").Value()

$b2 = $totalSheet.Range("B2").Value()
$c2 = $totalSheet.Range("C2").Value()
$d2 = $totalSheet.Range("D2").Value()

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = $b4
$totalSheet.Range("C5").Value = $c4
$totalSheet.Range("D5").Value = $d4

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = $b3
$totalSheet.Range("C4").Value = $c3
$totalSheet.Range("D4").Value = $d3

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = $b2
$totalSheet.Range("C3").Value = $c2
$totalSheet.Range("D3").Value = $d2

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.18

# Carry the "column A" style (bold/centered/bordered) down onto the newly
# created row 5.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)
